$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet="展览"; Cell="F5"; Old=2893; New=2898},
    @{Sheet="展览"; Cell="F7"; Old=233; New=234},
    @{Sheet="展览"; Cell="F9"; Old=291; New=292},
    @{Sheet="展览"; Cell="F10"; Old=6741; New=6753},
    @{Sheet="展览"; Cell="F11"; Old=27; New=28},
    @{Sheet="展览"; Cell="F12"; Old=0; New=4},
    @{Sheet="展览"; Cell="F13"; Old=318; New=321},
    @{Sheet="展览"; Cell="F14"; Old=585; New=588},
    @{Sheet="展览"; Cell="F15"; Old=1459; New=1465},
    @{Sheet="展览"; Cell="F17"; Old=1093; New=1095},
    @{Sheet="展览"; Cell="F18"; Old=2186; New=2193},
    @{Sheet="展览"; Cell="F19"; Old=1430; New=1437},
    @{Sheet="展览"; Cell="F20"; Old=639; New=640},
    @{Sheet="展览"; Cell="F21"; Old=90; New=91},
    @{Sheet="展览"; Cell="F22"; Old=1070; New=1074},
    @{Sheet="展览"; Cell="F23"; Old=76; New=78},
    @{Sheet="展览"; Cell="F24"; Old=154; New=156},
    @{Sheet="展览"; Cell="F25"; Old=309; New=311},
    @{Sheet="展览"; Cell="F26"; Old=1634; New=1639},
    @{Sheet="展览"; Cell="F27"; Old=1587; New=1617},
    @{Sheet="展览"; Cell="F28"; Old=532; New=533},
    @{Sheet="展览"; Cell="F30"; Old=27; New=28},
    @{Sheet="展览"; Cell="F31"; Old=1643; New=1644},
    @{Sheet="展览"; Cell="F32"; Old=1165; New=1173},
    @{Sheet="展览"; Cell="F34"; Old=572; New=573},
    @{Sheet="展览"; Cell="F35"; Old=10; New=12},
    @{Sheet="展览"; Cell="F37"; Old=381; New=385},
    @{Sheet="展览"; Cell="F38"; Old=2404; New=2410},
    @{Sheet="展览"; Cell="F39"; Old=2669; New=2676},
    @{Sheet="展览"; Cell="F41"; Old=172; New=173},
    @{Sheet="展览"; Cell="F44"; Old=14; New=15},
    @{Sheet="展览"; Cell="F46"; Old=114; New=115},
    @{Sheet="展览"; Cell="F47"; Old=153; New=154},
    @{Sheet="展览"; Cell="F48"; Old=128; New=130},
    @{Sheet="演出"; Cell="F7"; Old=130; New=131},
    @{Sheet="演出"; Cell="F10"; Old=29; New=30},
    @{Sheet="演出"; Cell="F12"; Old=174; New=175},
    @{Sheet="演出"; Cell="F14"; Old=53; New=54},
    @{Sheet="演出"; Cell="F15"; Old=50; New=52},
    @{Sheet="演出"; Cell="F20"; Old=11; New=13},
    @{Sheet="演出"; Cell="F23"; Old=448; New=449},
    @{Sheet="本地生活"; Cell="F6"; Old=1721; New=1722},
    @{Sheet="本地生活"; Cell="F7"; Old=1629; New=1637},
    @{Sheet="本地生活"; Cell="F8"; Old=1840; New=1839},
    @{Sheet="本地生活"; Cell="F9"; Old=2675; New=2679},
    @{Sheet="本地生活"; Cell="F10"; Old=974; New=977},
    @{Sheet="本地生活"; Cell="F11"; Old=865; New=871},
    @{Sheet="本地生活"; Cell="F12"; Old=34; New=35},
    @{Sheet="本地生活"; Cell="F13"; Old=212; New=216},
    @{Sheet="本地生活"; Cell="F14"; Old=550; New=1117},
    @{Sheet="本地生活"; Cell="F15"; Old=3703; New=6601},
    @{Sheet="全部类型"; Cell="F4"; Old=1721; New=1722},
    @{Sheet="全部类型"; Cell="F7"; Old=2893; New=2898},
    @{Sheet="全部类型"; Cell="F8"; Old=233; New=234},
    @{Sheet="全部类型"; Cell="F9"; Old=1629; New=1637},
    @{Sheet="全部类型"; Cell="F10"; Old=291; New=292},
    @{Sheet="全部类型"; Cell="F11"; Old=2675; New=2679},
    @{Sheet="全部类型"; Cell="F12"; Old=6741; New=6753},
    @{Sheet="全部类型"; Cell="F13"; Old=974; New=977},
    @{Sheet="全部类型"; Cell="F14"; Old=865; New=871},
    @{Sheet="全部类型"; Cell="F15"; Old=27; New=28},
    @{Sheet="全部类型"; Cell="F16"; Old=318; New=321},
    @{Sheet="全部类型"; Cell="F17"; Old=130; New=131},
    @{Sheet="全部类型"; Cell="F18"; Old=212; New=216},
    @{Sheet="全部类型"; Cell="F19"; Old=1093; New=1095},
    @{Sheet="全部类型"; Cell="F20"; Old=2186; New=2193},
    @{Sheet="全部类型"; Cell="F21"; Old=1430; New=1437},
    @{Sheet="全部类型"; Cell="F22"; Old=639; New=640},
    @{Sheet="全部类型"; Cell="F23"; Old=90; New=91},
    @{Sheet="全部类型"; Cell="F24"; Old=1070; New=1074},
    @{Sheet="全部类型"; Cell="F25"; Old=76; New=78},
    @{Sheet="全部类型"; Cell="F26"; Old=309; New=311},
    @{Sheet="全部类型"; Cell="F27"; Old=53; New=54},
    @{Sheet="全部类型"; Cell="F28"; Old=1634; New=1639},
    @{Sheet="全部类型"; Cell="F30"; Old=27; New=28},
    @{Sheet="全部类型"; Cell="F31"; Old=1643; New=1644},
    @{Sheet="全部类型"; Cell="F32"; Old=1165; New=1173},
    @{Sheet="全部类型"; Cell="F33"; Old=572; New=573},
    @{Sheet="全部类型"; Cell="F35"; Old=448; New=449},
    @{Sheet="全部类型"; Cell="F36"; Old=381; New=385},
    @{Sheet="全部类型"; Cell="F39"; Old=2404; New=2410},
    @{Sheet="全部类型"; Cell="F40"; Old=2669; New=2676},
    @{Sheet="全部类型"; Cell="F42"; Old=172; New=173},
    @{Sheet="全部类型"; Cell="F44"; Old=114; New=115},
    @{Sheet="全部类型"; Cell="F45"; Old=153; New=154}
)


foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.New
}

# Special case: sheet "本地生活" (Local Life) row 15 is now temporarily sold out,
# so its minimum price column G switches from a numeric price to status text.
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("G15").Value = "暂时售罄"
